# Update metric values across the existing sheets (1=Bidirectional A, 2=D Lite, 3=IDA, 4=SMA, 7=RTAA (L=25, M=3))
$wb = $excel.ActiveWorkbook

# Sheet 1: Bidirectional A
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = 0.0001506669996160781
$ws.Range("B3").Value = 0.0001620410002942663
$ws.Range("B4").Value = 2677.5408
$ws.Range("B5").Value = 2790.3315
$ws.Range("B6").Value = 0.00287628173828125
$ws.Range("B7").Value = 0.00285797119140625
$ws.Range("B8").Value = 0.0084381103515625
$ws.Range("B9").Value = 0.0084381103515625
$ws.Range("B10").Value = 0.00004329200055508409
$ws.Range("B11").Value = 0.0000559170002816245
$ws.Range("B12").Value = 0.00004169159983575809
$ws.Range("B13").Value = 0.00004960019978170749
$ws.Range("B14").Value = 5

# Sheet 2: D Lite
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = 0.0001537499992991798
$ws.Range("B3").Value = 0.001023999999233638
$ws.Range("B4").Value = 2677.5408
$ws.Range("B5").Value = 2677.5408
$ws.Range("B6").Value = 0.0028533935546875
$ws.Range("B7").Value = 0.0028533935546875
$ws.Range("B8").Value = 0.07726287841796875
$ws.Range("B9").Value = 0.07698516845703125
$ws.Range("B10").Value = 0.0002467909998813411
$ws.Range("B11").Value = 0.0008286670017696451
$ws.Range("B12").Value = 0.00004184159988653846
$ws.Range("B13").Value = 0.0000428502000431763
$ws.Range("B14").Value = 5

# Sheet 3: IDA
$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = 0.00025558400011505
$ws.Range("B3").Value = 0.01366037499974482
$ws.Range("B4").Value = 2677.5408
$ws.Range("B5").Value = 2677.5408
$ws.Range("B6").Value = 0.0028533935546875
$ws.Range("B7").Value = 0.0028533935546875
$ws.Range("B8").Value = 0.0026397705078125
$ws.Range("B9").Value = 0.0025848388671875
$ws.Range("B10").Value = 0.00004600000102072954
$ws.Range("B11").Value = 0.01290912500007835
$ws.Range("B12").Value = 0.01184026660012023
$ws.Range("B13").Value = 0.0000402916000894038
$ws.Range("B14").Value = 5

# Sheet 4: SMA
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = 0.0001817499996832339
$ws.Range("B3").Value = 0.00007399999958579428
$ws.Range("B4").Value = 2677.5408
$ws.Range("B5").Value = 2677.5408
$ws.Range("B6").Value = 0.00287628173828125
$ws.Range("B7").Value = 0.00285797119140625
$ws.Range("B10").Value = 0.00004366699977254029
$ws.Range("B11").Value = 0.00005908300045120995
$ws.Range("B12").Value = 0.00005871680041309446
$ws.Range("B13").Value = 0.00003793340038100723
$ws.Range("B14").Value = 5

# Sheet 7: RTAA (L=25, M=3)
$ws = $wb.Worksheets.Item(7)
$ws.Range("B2").Value = 0.0001637500008655479
$ws.Range("B3").Value = 0.000217417000385467
$ws.Range("B4").Value = 2677.5408
$ws.Range("B5").Value = 3484.7703
$ws.Range("B6").Value = 0.00290679931640625
$ws.Range("B7").Value = 0.00286407470703125
$ws.Range("B8").Value = 0.00658416748046875
$ws.Range("B9").Value = 0.00658416748046875
$ws.Range("B10").Value = 0.00004354100019554608
$ws.Range("B11").Value = 0.0001148749997810228
$ws.Range("B12").Value = 0.0001135832004365511
$ws.Range("B13").Value = 0.00003740840002137702
$ws.Range("B14").Value = 5

# IDA sheet column B width: 24 -> 23 (raw OOXML width units)
$ws = $wb.Worksheets.Item(3)
$ws.Columns.Item(2).ColumnWidth = 22.17

# Add new summary sheet "A" at the end of the workbook
$template = $wb.Worksheets.Item(1)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "A"

# Copy header/data formatting from the template sheet so styles match (bold centered header, centered data)
$template.Range("A1:B1").Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122)
$template.Range("A2:B2").Copy()
$newSheet.Range("A2:B3").PasteSpecial(-4122)

$newSheet.Range("A1").Value = "Metric"
$newSheet.Range("B1").Value = "Value"
$newSheet.Range("A2").Value = "Cost"
$newSheet.Range("B2").Value = 2677.5408
$newSheet.Range("A3").Value = "Path length"
$newSheet.Range("B3").Value = 5

$newSheet.Columns.Item(1).ColumnWidth = 12.17
$newSheet.Columns.Item(2).ColumnWidth = 19.17

$newSheet.Range("A1").Select()
